# Generate Report for Archive
#
# The localization status text "Ready for handoff" has moved on to
# "In Translation" for both tracked files/languages. Update every cell
# that shows this status (the Overview sheet's per-language status
# columns, plus each language sheet's own "Status" column) and re-fit
# the status column widths to the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status is reported per-language in columns E (zh-cn)
# and F (de-de) for each of the two tracked files (rows 2-3). ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $overview.Range($addr)
    if ($cell.Text -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- Per-language sheets: column C holds the "Status" value for each
# tracked file (rows 2-3). ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3")) {
        $cell = $ws.Range($addr)
        if ($cell.Text -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Re-fit the status columns now that the text is shorter. The
# runtime's ColumnWidth setter snaps to the nearest display-pixel
# increment, so feed it the character width that lands closest to the
# desired ~13.41 "characters" target width. ---
$statusColWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $statusColWidth  # E: zh-cn status
$overview.Columns.Item(6).ColumnWidth = $statusColWidth  # F: de-de status

$wb.Worksheets.Item("zh-cn").Columns.Item(3).ColumnWidth = $statusColWidth
$wb.Worksheets.Item("de-de").Columns.Item(3).ColumnWidth = $statusColWidth
